$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.308.37'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.863.29'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''0.7015'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '''238.05'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '''0.08217'
$ws.Range("E8").Value = '  +9.82%  '
$ws.Range("D9").Value = '''0.3049'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '''23.35'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").Value = '''0.08173'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '1.875.55'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '''0.7187'
$ws.Range("D14").Value = '''5.184'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").Value = '''89.43'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '29.320.95'
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.000007891'
$ws.Range("E17").Value = '  +2.86%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '''5.784'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '''13.42'
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("D20").Value = '''237.79'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = '''0.9999'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '2.106.15'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '''7.471'
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '''162.35'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").Value = '''9.001'
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = '''0.1444'
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").Value = '''1.991'
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("D30").Value = '''1.435'
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("D31").Value = '''4.436'
$ws.Range("E31").Value = '  -2.76%  '
$ws.Range("D32").Value = '''1.487'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").Value = '''4.063'
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").Value = '''0.7059'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '''1.006'
$ws.Range("E37").Value = '  -2.44%  '
$ws.Range("D38").Value = '''2.662'
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = '''0.01856'
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").Value = '''2.720'
$ws.Range("E40").Value = '  +1.64%  '
$ws.Range("D41").Value = '1.145.80'
$ws.Range("E41").Value = '  +7.76%  '
$ws.Range("D42").Value = '''0.9222'
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").Value = '''5.974'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '''0.4284'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '''70.99'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").Value = '''0.9999'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '''103.06'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '''1.778'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").Value = '2.004.53'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("D50").Value = '''9.218'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '''6.986'
$ws.Range("E51").Value = '  -1.09%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
